# Add FEBRERO week-7 data rows (169-195) to Sheet1, mirroring the
# week-6 block directly above it, and move the viewport/selection to
# reflect the newly-entered data (matches Excel's behaviour after
# typing a block of rows and leaving the cursor below it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$weekData = @(
    @(169, 2026, "FEBRERO", 7, "AMARILLO", "COLORES", "GOLDFINCH", 6018),
    @(170, 2026, "FEBRERO", 7, "AMARILLO", "COLORES", "HIGH AND EXOTIC", 14570),
    @(171, 2026, "FEBRERO", 7, "AMARILLO", "COLORES", "MOMENTUM", 8570),
    @(172, 2026, "FEBRERO", 7, "AMARILLO", "COLORES", "SUNDAY MORNING", 342),
    @(173, 2026, "FEBRERO", 7, "BIC. AMARILLO", "COLORES", "SUMMER LIGHT", 2188),
    @(174, 2026, "FEBRERO", 7, "BICOLOR", "COLORES", "BLUSH", 1430),
    @(175, 2026, "FEBRERO", 7, "BICOLOR", "COLORES", "DISCOVERY", 1050),
    @(176, 2026, "FEBRERO", 7, "BLANCO", "COLORES", "HIGH AND PURE", 19920),
    @(177, 2026, "FEBRERO", 7, "BLANCO", "COLORES", "SUGGAR DOLL", 4630),
    @(178, 2026, "FEBRERO", 7, "BLANCO", "COLORES", "VANILLA ICE", 9350),
    @(179, 2026, "FEBRERO", 7, "BLANCO", "COLORES", "VENDELA", 4507),
    @(180, 2026, "FEBRERO", 7, "DURAZNO", "COLORES", "TIFANY", 2350),
    @(181, 2026, "FEBRERO", 7, "HOT PINK", "COLORES", "COTTON CANDY", 4145),
    @(182, 2026, "FEBRERO", 7, "HOT PINK", "COLORES", "JACARANDA", 10863),
    @(183, 2026, "FEBRERO", 7, "HOT PINK", "COLORES", "PINK FLOYD", 14371),
    @(184, 2026, "FEBRERO", 7, "LAVANDER", "COLORES", "DEEP PURPLE", 17700),
    @(185, 2026, "FEBRERO", 7, "LAVANDER", "COLORES", "MOODY BLUES", 7545),
    @(186, 2026, "FEBRERO", 7, "NARANJA", "COLORES", "ALIVE", 1480),
    @(187, 2026, "FEBRERO", 7, "NARANJA", "COLORES", "BROMO", 3100),
    @(188, 2026, "FEBRERO", 7, "NARANJA", "COLORES", "CLEMENTINA", 2947),
    @(189, 2026, "FEBRERO", 7, "NARANJA", "COLORES", "NINA", 26178),
    @(190, 2026, "FEBRERO", 7, "ROJO", "ROJO", "FREEDOM", 123902),
    @(191, 2026, "FEBRERO", 7, "ROSADO", "COLORES", "ABSOLUT IN PINK", 3600),
    @(192, 2026, "FEBRERO", 7, "ROSADO", "COLORES", "HIGH AND BONITA", 11948),
    @(193, 2026, "FEBRERO", 7, "ROSADO", "COLORES", "LUCIANO", 2160),
    @(194, 2026, "FEBRERO", 7, "ROSADO", "COLORES", "STARFISH", 5658),
    @(195, 2026, "FEBRERO", 7, "ROSADO", "COLORES", "TABATHA", 9021)
)

foreach ($r in $weekData) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
}

# Move the viewport / selection the way Excel would after this edit.
$excel.ActiveWindow.ScrollRow = 160
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A196").Select()
